$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newTexts = @(
    "70-64=",
    "92-37=",
    "33-25=",
    "16+56=",
    "92-37=",
    "36+45=",
    "22+69=",
    "80-45=",
    "60-4=",
    "24-19=",
    "55+38=",
    "22-16=",
    "90-53=",
    "12+19=",
    "74+19=",
    "71-63=",
    "36+47=",
    "18+15=",
    "55-17=",
    "19+24=",
    "43+48=",
    "51-28=",
    "24+67=",
    "74-47=",
    "32+49=",
    "82-3=",
    "93-4=",
    "19+13=",
    "85-48=",
    "80-45=",
    "38+3=",
    "59+37=",
    "37+29=",
    "91-57=",
    "61-27=",
    "16+55=",
    "35+46=",
    "26+36=",
    "14+57=",
    "23+69=",
    "85-36=",
    "3+19=",
    "84-78=",
    "26+27=",
    "25-16=",
    "45+7=",
    "72-46=",
    "65+16=",
    "54+27=",
    "40-11=",
    "11-4=",
    "32-14=",
    "61-49=",
    "88-79=",
    "29+47=",
    "2+59=",
    "40-12=",
    "24-17=",
    "90-88=",
    "67-19=",
    "85-37=",
    "60-53=",
    "31-6=",
    "29+58=",
    "65-6=",
    "92-25=",
    "7+48=",
    "47+4=",
    "39+54=",
    "27+5=",
    "45-39=",
    "43+29=",
    "82-48=",
    "13+58=",
    "57+25=",
    "95-86=",
    "61-2=",
    "37+19=",
    "88-49=",
    "24+7=",
    "24+67=",
    "64-29=",
    "62-38=",
    "33-8=",
    "50-27=",
    "54-35=",
    "68+8=",
    "62-6=",
    "29+27=",
    "74-39=",
    "45+29=",
    "77+19=",
    "45-29=",
    "64-38=",
    "48+24=",
    "80-67=",
    "76-27=",
    "85-26=",
    "26+26=",
    "26+27="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newTexts[$idx]
        $idx = $idx + 1
    }
}

Write-Output "done: $idx cells updated"
